$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the "remaining" counter (column E) simply decrements by 1
# (daily countdown tick; start date in column F is unchanged).
$simpleDecrement = @{
    2 = 11
    3 = 11
    4 = 11
    5 = 7
    6 = 11
    7 = 7
    8 = 11
    9 = 7
    10 = 4
    11 = 11
    12 = 7
    13 = 11
    14 = 11
    15 = 11
    16 = 1
    17 = 7
    22 = 7
    23 = 7
    24 = 7
    25 = 7
    26 = 7
    27 = 5
    40 = 4
    41 = 4
    43 = 7
    44 = 4
    45 = 7
    46 = 4
    48 = 4
    49 = 5
    50 = 5
    51 = 5
    52 = 5
    53 = 5
    54 = 5
    55 = 5
    56 = 5
    57 = 5
    58 = 9
    59 = 9
    60 = 9
    61 = 5
    62 = 9
    63 = 9
    64 = 9
    70 = 1
    71 = 1
    72 = 1
    73 = 1
    74 = 1
    75 = 1
    76 = 1
    77 = 4
    78 = 4
    79 = 4
    80 = 4
    81 = 4
    82 = 4
    83 = 4
    84 = 4
    85 = 4
    86 = 4
    87 = 4
    88 = 4
    89 = 4
    90 = 4
    91 = 7
    92 = 4
    93 = 4
    95 = 3
    96 = 1
    97 = 1
    98 = 1
    99 = 1
}

foreach ($row in $simpleDecrement.Keys) {
    $ws.Cells.Item($row, 5).Value = $simpleDecrement[$row]
}

# Rows where the counter hit its floor and got renewed: remaining (E) resets
# back up to the total (D), and the start date (F) rolls to the new cycle start.
$renewals = @{
    18 = @(10, 20251218)
    19 = @(10, 20251218)
    20 = @(10, 20251218)
    21 = @(10, 20251218)
    28 = @(10, 20251218)
    29 = @(10, 20251218)
    30 = @(10, 20251218)
    31 = @(10, 20251218)
    32 = @(10, 20251218)
    33 = @(10, 20251218)
    34 = @(10, 20251218)
    35 = @(10, 20251218)
    37 = @(10, 20251218)
    38 = @(10, 20251218)
    39 = @(10, 20251218)
    42 = @(10, 20251218)
    47 = @(10, 20251218)
    65 = @(10, 20251218)
    66 = @(10, 20251218)
    67 = @(10, 20251218)
    68 = @(10, 20251218)
    69 = @(10, 20251218)
    94 = @(7, 20251218)
}

foreach ($row in $renewals.Keys) {
    $values = $renewals[$row]
    $ws.Cells.Item($row, 5).Value = $values[0]
    $ws.Cells.Item($row, 6).Value = $values[1]
}
